$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 is the only data row in the table; update it with the latest
# automatically-fetched electricity spot prices.
$ws.Range("A2").Value = 46049
$ws.Range("B2").Value = 18.46
$ws.Range("C2").Value = 5.12
$ws.Range("D2").Value = 3.78
$ws.Range("G2").Value = 3.78
$ws.Range("H2").Value = 4.84
$ws.Range("I2").Value = 17.63
$ws.Range("J2").Value = 32.77
$ws.Range("K2").Value = 33.14
$ws.Range("L2").Value = 26.63
$ws.Range("M2").Value = 26.88
$ws.Range("N2").Value = 26.34
$ws.Range("O2").Value = 22.72
$ws.Range("P2").Value = 20.23
$ws.Range("Q2").Value = 13.71
$ws.Range("R2").Value = 23.84
$ws.Range("S2").Value = 46.51
$ws.Range("T2").Value = 55.77
$ws.Range("U2").Value = 85.34
$ws.Range("V2").Value = 95.92
$ws.Range("W2").Value = 69.63
$ws.Range("X2").Value = 37.44
$ws.Range("Y2").Value = 13.59
$ws.Range("Z2").Value = 28.82
$ws.Range("AB2").Value = 54.14
$ws.Range("AD2").Value = 82.78
$ws.Range("AF2").Value = 70.56
